$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 2nd Sep typing test score update (row 14, "2nd Sep"):
# fill in Typing Test (WPM) and Accuracy (%), previously "-" placeholders.
$ws.Range("J14").Value = 44.3
$ws.Range("K14").Value = 96.43
